$d = $word.ActiveDocument

# 1. Activation date
$d.Content.Find.Execute(
    "Ativação: 01/01/2018", $true, $true, $false, $false, $false,
    $true, 1, $false, "Ativação: 01/01/2025", 2)

# 2. Portuguese summary paragraph
$d.Content.Find.Execute(
    "Espaços vetoriais, Transformações lineares, auto-valores e auto-vetores, diagonalização de Operadores, espaços vetoriais com produto interno, aplicações as equações diferenciais.",
    $true, $true, $false, $false, $false,
    $true, 1, $false,
    "Espaços vetoriais, transformações lineares, diagonalização de operadores lineares, forma canônica de Jordan, espaços vetoriais com produto interno, aplicações a sistemas dinâmicos.",
    2)

# 3. English (italic) summary paragraph
$d.Content.Find.Execute(
    "Vector Spaces, linear Transformations, eigenvalues and eigenvectors, Diagonalization, Inner product in vectorial Spaces, applications to differential Equations.",
    $true, $true, $false, $false, $false,
    $true, 1, $false,
    "Vector spaces, linear transformations, diagonalization of linear operators, Jordan canonical form, inner product spaces, applications to dynamical systems.",
    2)

# 4. Capitalize "spaços" -> "Espaços" at start of Programa paragraph
$d.Content.Find.Execute(
    "spaços vetoriais: Definição, Propriedades dos Espaços Vetoriais",
    $true, $true, $false, $false, $false,
    $true, 1, $false,
    "Espaços vetoriais: Definição, Propriedades dos Espaços Vetoriais",
    2)

# 5. Bibliography paragraph.
# The long <w:t> is rewritten entirely, but Word's Find text is capped at
# ~255 chars, so the change is applied as a sequence of short, unique
# anchored replacements that together produce the same result. Note: the
# stored text uses U+037E (GREEK QUESTION MARK, a semicolon look-alike) in
# a few spots; $gq below matches that exact character so Find locates it.
$gq = [string][char]0x037E

$d.Content.Find.Execute(
    "01.STRANG", $true, $true, $false, $false, $false,
    $true, 1, $false,
    "1.LIMA, Elon Lages. Álgebra Linear, IMPA, 2020. ISBN: 978-65-990528-3-5. 10ª edição.2.HOFFMAN, Kennethe; KUNZE, Ray. Linear Algebra. Pearson. 1971. 2nd Edition.3. STRANG",
    2)

$d.Content.Find.Execute(
    "2010.02.LIPSCHUTZ, Seymour . Algebra linear.", $true, $true, $false, $false, $false,
    $true, 1, $false,
    "2010.4.LIPSCHUTZ, Seymour. Álgebra linear.",
    2)

$d.Content.Find.Execute(
    ("1990.03.HOWARD, Anton " + $gq + " RORRES,Chris."), $true, $true, $false, $false, $false,
    $true, 1, $false,
    "1990.5.HOWARD, Anton ; RORRES, Chris.",
    2)

$d.Content.Find.Execute(
    "2001.04.MICHOLSON, W. Keith.", $true, $true, $false, $false, $false,
    $true, 1, $false,
    "2001.6.MICHOLSON, W. Keith.",
    2)

$d.Content.Find.Execute(
    "Mc GrawHill,  2006.05.BOLDRINI,", $true, $true, $false, $false, $false,
    $true, 1, $false,
    "Mc GrawHill, 2006.7.BOLDRINI,",
    2)

$d.Content.Find.Execute(
    ("BOLDRINI, José Luiz " + $gq + " COSTA Sueli I. Rodrigues" + $gq + " FIGUEIREDO Vera Lúcia" + $gq + " WETZLER"),
    $true, $true, $false, $false, $false,
    $true, 1, $false,
    "BOLDRINI, José Luiz ; COSTA Sueli I. Rodrigues; FIGUEIREDO Vera Lúcia; WETZLER",
    2)

$d.Content.Find.Execute(
    "1986.06.POOLE, David.", $true, $true, $false, $false, $false,
    $true, 1, $false,
    "1986.8.POOLE, David.",
    2)
